# acta_entregachip.docx edits
#  1. "LINEA ENTREGADA Nº: 0959659750" paragraph -> "LINEA ENTREGADA Nº:"
#     (drop the serial-number runs, split "Nº" into its own spell-checked
#     run, move the _GoBack bookmark here)
#  2. "CODIGO Nº:[...]" paragraph -> add gramStart/gramEnd proofErr markers
#     around the "Nº:...[" span and drop the _GoBack bookmark (it moved to
#     the previous paragraph)
#  3. Update the printed date field result 02/12/2024 -> 22/01/2025

$d = $word.ActiveDocument
$nord = [char]0x00BA   # 'º' MASCULINE ORDINAL INDICATOR

# ---------------------------------------------------------------------
# 1. "LINEA ENTREGADA Nº:" paragraph
# ---------------------------------------------------------------------
$find1 = $d.Content.Duplicate
$found1 = $find1.Find.Execute("LINEA ENTREGADA N" + $nord + ":*0959659750", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the 'LINEA ENTREGADA' paragraph"
}
# Extend by one character to also capture the paragraph mark, so
# InsertXML replaces the whole paragraph (incl. its pPr) cleanly.
$r1 = $d.Range($find1.Start, $find1.End + 1)

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D14F73" w:rsidRPr="005A773C" w:rsidRDefault="00D14F73" w:rsidP="00A53633">' +
        '<w:pPr><w:tabs><w:tab w:val="left" w:pos="1200"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/>' +
        '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t xml:space="preserve">LINEA ENTREGADA </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>N' + $nord + '</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>:</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
        '</w:p>'

$r1.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2. "CODIGO Nº:[cuadrilla.equipos]" paragraph
# ---------------------------------------------------------------------
$find2 = $d.Content.Duplicate
$found2 = $find2.Find.Execute("CODIGO N" + $nord + ":*cuadrilla.*equipos*", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the 'CODIGO' paragraph"
}
# Extend by one character to also capture the paragraph mark, so
# InsertXML replaces the whole paragraph (incl. its pPr) cleanly.
$r2 = $d.Range($find2.Start, $find2.End + 1)

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D14F73" w:rsidRDefault="00D14F73" w:rsidP="003B2B5F">' +
        '<w:pPr><w:tabs><w:tab w:val="left" w:pos="1200"/></w:tabs><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/>' +
        '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t xml:space="preserve">CODIGO </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>N' + $nord + '</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r w:rsidRPr="00A53633"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>:</w:t></w:r>' +
        '<w:r w:rsidR="00B83774" w:rsidRPr="004C1157"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>[</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/>' +
        '<w:r w:rsidR="00B83774" w:rsidRPr="004C1157"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>cuadrilla.</w:t></w:r>' +
        '<w:r w:rsidR="00B83774"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>equipos</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r w:rsidR="00B83774" w:rsidRPr="004C1157"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="24"/></w:rPr>' +
        '<w:t>]</w:t></w:r>' +
        '</w:p>'

$r2.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3. Date field result: 02/12/2024 -> 22/01/2025
# ---------------------------------------------------------------------
$d.Content.Find.Execute("02/12/2024", $false, $false, $false, $false, $false, $true, 1, $false, "22/01/2025", 2)
